$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.466.93'
$ws.Range("E2").Value = '  +2.03%  '

# Row 3
$ws.Range("D3").Value = '1.838.37'
$ws.Range("E3").Value = '  +0.92%  '

# Row 4
$ws.Range("E4").Value = '  -0.16%  '

# Row 5
$ws.Range("D5").Value = '''243.55'
$ws.Range("E5").Value = '  +1.13%  '

# Row 6
$ws.Range("D6").Value = '''0.6236'
$ws.Range("E6").Value = '  +1.41%  '

# Row 7
$ws.Range("D7").Value = '''0.9991'
$ws.Range("E7").Value = '  +0.34%  '

# Row 8
$ws.Range("D8").Value = '''0.07396'
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("D9").Value = '''0.2933'
$ws.Range("E9").Value = '  +0.51%  '

# Row 10
$ws.Range("D10").Value = '''23.44'
$ws.Range("E10").Value = '  +2.35%  '

# Row 11
$ws.Range("D11").Value = '''0.07659'
$ws.Range("E11").Value = '  +0.50%  '

# Row 12
$ws.Range("D12").Value = '1.846.78'
$ws.Range("E12").Value = '  +1.51%  '

# Row 13
$ws.Range("D13").Value = '''5.013'
$ws.Range("E13").Value = '  +0.96%  '

# Row 14
$ws.Range("D14").Value = '''0.6798'
$ws.Range("E14").Value = '  +1.70%  '

# Row 15
$ws.Range("D15").Value = '''83.17'
$ws.Range("E15").Value = '  +0.94%  '

# Row 16
$ws.Range("D16").Value = '''0.000008968'
$ws.Range("E16").Value = '  -0.53%  '

# Row 17
$ws.Range("D17").Value = '''5.915'
$ws.Range("E17").Value = '  +1.07%  '

# Row 18
$ws.Range("D18").Value = '29.440.49'
$ws.Range("E18").Value = '  +1.90%  '

# Row 19
$ws.Range("D19").Value = '2.092.66'
$ws.Range("E19").Value = '  -1.58%  '

# Row 20
$ws.Range("D20").Value = '''244.57'
$ws.Range("E20").Value = '  +1.43%  '

# Row 21
$ws.Range("D21").Value = '''12.57'
$ws.Range("E21").Value = '  -0.40%  '

# Row 22
$ws.Range("D22").Value = '''0.9995'
$ws.Range("E22").Value = '  +0.39%  '

# Row 23
$ws.Range("D23").Value = '''7.413'
$ws.Range("E23").Value = '  +3.25%  '

# Row 24
$ws.Range("D24").Value = '''0.9992'
$ws.Range("E24").Value = '  +0.07%  '

# Row 25
$ws.Range("D25").Value = '''158.54'
$ws.Range("E25").Value = '  +0.54%  '

# Row 26
$ws.Range("B26").Value = 'Stellar'
$ws.Range("C26").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D26").Value = '''0.1400'
$ws.Range("E26").Value = '  -0.59%  '

# Row 27
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").Value = '''8.587'
$ws.Range("E27").Value = '  +1.81%  '

# Row 28
$ws.Range("D28").Value = '''17.78'
$ws.Range("E28").Value = '  +0.13%  '

# Row 29
$ws.Range("D29").Value = '''1.496'
$ws.Range("E29").Value = '  +0.97%  '

# Row 30
$ws.Range("D30").Value = '''0.05911'
$ws.Range("E30").Value = '  +6.26%  '

# Row 31
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").Value = '''1.230'
$ws.Range("E31").Value = '  +2.07%  '

# Row 32
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '''4.102'
$ws.Range("E32").Value = '  +0.47%  '

# Row 33
$ws.Range("B33").Value = 'Filecoin'
$ws.Range("C33").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D33").Value = '''4.118'
$ws.Range("E33").Value = '  +0.45%  '

# Row 34
$ws.Range("E34").Value = '  +2.38%  '

# Row 35
$ws.Range("D35").Value = '''1.143'
$ws.Range("E35").Value = '  +0.96%  '

# Row 36
$ws.Range("D36").Value = '''0.7245'
$ws.Range("E36").Value = '  -1.63%  '

# Row 37
$ws.Range("D37").Value = '''2.611'
$ws.Range("E37").Value = '  -0.43%  '

# Row 38
$ws.Range("D38").Value = '''2.873'
$ws.Range("E38").Value = '  +4.38%  '

# Row 39
$ws.Range("D39").Value = '1.225.36'
$ws.Range("E39").Value = '  +2.00%  '

# Row 40
$ws.Range("D40").Value = '''0.01767'
$ws.Range("E40").Value = '  -0.08%  '

# Row 41
$ws.Range("D41").Value = '''0.9170'
$ws.Range("E41").Value = '  +2.72%  '

# Row 42
$ws.Range("D42").Value = '''6.244'
$ws.Range("E42").Value = '  -1.58%  '

# Row 43
$ws.Range("E43").Value = '  +0.64%  '

# Row 44
$ws.Range("D44").Value = '2.006.46'
$ws.Range("E44").Value = '  +0.04%  '

# Row 45
$ws.Range("D45").Value = '''101.98'
$ws.Range("E45").Value = '  +1.08%  '

# Row 46
$ws.Range("D46").Value = '''65.82'
$ws.Range("E46").Value = '  +1.40%  '

# Row 47
$ws.Range("D47").Value = '''0.00000000122'
$ws.Range("E47").Value = '  +4.37%  '

# Row 48
$ws.Range("D48").Value = '''0.5056'
$ws.Range("E48").Value = '  -0.06%  '

# Row 49
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").Value = '''0.4062'
$ws.Range("E49").Value = '  +0.68%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").Value = '''9.200'
$ws.Range("E50").Value = '  +1.54%  '

# Row 51
$ws.Range("D51").Value = '''0.1172'
$ws.Range("E51").Value = '  +7.00%  '
